$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update selection on sheet1 (Blad1) before adding the new sheet
$ws1.Range("N10").Select()

# Add the new worksheet "Blad2" after Blad1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Blad2"

# Column widths
$ws2.Columns.Item(1).ColumnWidth = 14.5703125
$ws2.Columns.Item(2).ColumnWidth = 18.28515625
$ws2.Columns.Item(3).ColumnWidth = 27.42578125
$ws2.Columns.Item(4).ColumnWidth = 14.5703125
$ws2.Columns.Item(5).ColumnWidth = 9.85546875
$ws2.Columns.Item(6).ColumnWidth = 11.85546875
$ws2.Columns.Item(7).ColumnWidth = 71.5703125

# Title
$ws2.Range("A1").Value = 'Backlog items'
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").Font.Size = 12

# Header row (row 3) - shared strings 37-42, styled with built-in "Good" style
$ws2.Range("B3").Value = 'Item number'
$ws2.Range("C3").Value = 'Item'
$ws2.Range("D3").Value = 'Estimated time'
$ws2.Range("E3").Value = 'Priority'
$ws2.Range("F3").Value = 'Status'
$ws2.Range("G3").Value = 'Description'
$ws2.Range("A3:G3").Style = "Good"

# Backlog rows 4-9: write F (Status) and G (Description) first for row 4 to create
# the shared-string "ToDo" and its first description ahead of the Item name, then
# continue left-to-right for the remaining columns (matches original authoring order).
$ws2.Range("B4").Value = 1
$ws2.Range("F4").Value = 'ToDo'
$ws2.Range("G4").Value = 'Köra med en test implementering av Identity med Google'
$ws2.Range("C4").Value = 'Test Identity'
$ws2.Range("D4").Value = 2
$ws2.Range("E4").Value = 1

$ws2.Range("B5").Value = 2
$ws2.Range("C5").Value = 'Implementera Identity'
$ws2.Range("G5").Value = 'Implementera Identity, få igång alla funktioner'
$ws2.Range("F5").Value = 'ToDo'
$ws2.Range("D5").Value = 4
$ws2.Range("E5").Value = 1

$ws2.Range("B6").Value = 3
$ws2.Range("C6").Value = 'Modulera Databas'
$ws2.Range("G6").Value = 'Finslipa, modulera slutgiltlig databas'
$ws2.Range("F6").Value = 'ToDo'
$ws2.Range("D6").Value = 2
$ws2.Range("E6").Value = 1

$ws2.Range("B7").Value = 4
$ws2.Range("C7").Value = 'Implementera rest databas'
$ws2.Range("G7").Value = 'Implementering av resterande databas med entity framework code first'
$ws2.Range("F7").Value = 'ToDo'
$ws2.Range("D7").Value = 2
$ws2.Range("E7").Value = 1

$ws2.Range("B8").Value = 5
$ws2.Range("C8").Value = 'Layout login-sida'
$ws2.Range("G8").Value = 'Styla login sidan med Material Design'
$ws2.Range("F8").Value = 'ToDo'
$ws2.Range("D8").Value = 4
$ws2.Range("E8").Value = 1

$ws2.Range("B9").Value = 6
$ws2.Range("C9").Value = 'Implementera wep-api'
$ws2.Range("G9").Value = 'Implementer api/service för alla klasser/tabeller, och dess respektive angular'
$ws2.Range("F9").Value = 'ToDo'
$ws2.Range("D9").Value = 6
$ws2.Range("E9").Value = 1

# Backlog rows 11-15
$ws2.Range("B11").Value = 8
$ws2.Range("C11").Value = 'Implementering av MaterDes'
$ws2.Range("G11").Value = 'Implementer resterande design, Material Design'
$ws2.Range("F11").Value = 'ToDo'
$ws2.Range("D11").Value = 8
$ws2.Range("E11").Value = 1

$ws2.Range("B12").Value = 9
$ws2.Range("C12").Value = 'Skriva testcases'
$ws2.Range("G12").Value = 'Skriva testcases för sidan'
$ws2.Range("F12").Value = 'ToDo'
$ws2.Range("D12").Value = 2
$ws2.Range("E12").Value = 2

$ws2.Range("B13").Value = 10
$ws2.Range("C13").Value = 'Köra igenom testcases'
$ws2.Range("G13").Value = 'Köra igenom alla testcases för sidan'
$ws2.Range("F13").Value = 'ToDo'
$ws2.Range("D13").Value = 2
$ws2.Range("E13").Value = 2

$ws2.Range("B14").Value = 11
$ws2.Range("C14").Value = 'Fixa till eventuella fel'
$ws2.Range("G14").Value = 'Rätta till alla eventuella fel som vi har upptäckt'
$ws2.Range("F14").Value = 'ToDo'
$ws2.Range("D14").Value = 2
$ws2.Range("E14").Value = 2

$ws2.Range("B15").Value = 12
$ws2.Range("C15").Value = 'Implenetera loggning av fel'
$ws2.Range("G15").Value = 'Implementera funktion som loggar alla fel och exceptions'
$ws2.Range("F15").Value = 'ToDo'
$ws2.Range("D15").Value = 1
$ws2.Range("E15").Value = 3

# Row 10 (inserted last, after the others - matches shared string indices 66/67)
$ws2.Range("B10").Value = 7
$ws2.Range("C10").Value = 'Implementera div funktioner'
$ws2.Range("G10").Value = 'Implementering av diverse funktioner, features'
$ws2.Range("F10").Value = 'ToDo'
$ws2.Range("D10").Value = 4
$ws2.Range("E10").Value = 1

# Selection on the new sheet, which becomes the active tab
$ws2.Range("G10").Select()

Write-Host "done"
